$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '63.374.21'
$ws.Range('E2').Value = '  -1.54%  '

# Row 3
$ws.Range('D3').Value = '2.576.40'
$ws.Range('E3').Value = '  -3.00%  '

# Row 4
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').Value = '587.14'
$ws.Range('E5').Value = '  -3.82%  '

# Row 6
$ws.Range('D6').Value = '149.12'
$ws.Range('E6').Value = '  -1.32%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('E8').Value = '  -0.96%  '

# Row 9
$ws.Range('E9').Value = '  -0.37%  '

# Row 10
$ws.Range('D10').Value = '5.75'
$ws.Range('E10').Value = '  +2.29%  '

# Row 11
$ws.Range('E11').Value = '  -1.81%  '

# Row 12
$ws.Range('E12').Value = '  -0.84%  '

# Row 13
$ws.Range('D13').Value = '27.47'
$ws.Range('E13').Value = '  -1.65%  '

# Row 14
$ws.Range('D14').Value = '3.039.87'
$ws.Range('E14').Value = '  -2.97%  '

# Row 15
$ws.Range('D15').Value = '63.237.38'
$ws.Range('E15').Value = '  -1.51%  '

# Row 16
$ws.Range('E16').Value = '  +3.93%  '

# Row 17
$ws.Range('D17').Value = '2.579.19'
$ws.Range('E17').Value = '  -2.62%  '

# Row 18
$ws.Range('D18').Value = '12.16'
$ws.Range('E18').Value = '  +0.44%  '

# Row 19
$ws.Range('D19').Value = '4.67'
$ws.Range('E19').Value = '  +0.78%  '

# Row 20
$ws.Range('D20').Value = '345.12'
$ws.Range('E20').Value = '  -0.47%  '

# Row 21
$ws.Range('E21').Value = '  -1.58%  '

# Row 22
$ws.Range('E22').Value = '  -0.03%  '

# Row 23
$ws.Range('E23').Value = '  +0.18%  '

# Row 24
$ws.Range('D24').Value = '1.69'
$ws.Range('E24').Value = '  -5.26%  '

# Row 25
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = '9.10'
$ws.Range('E25').Value = '  -3.37%  '

# Row 26
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').Value = '1.66'
$ws.Range('E26').Value = '  -3.78%  '

# Row 27
$ws.Range('D27').Value = '553.40'
$ws.Range('E27').Value = '  -0.51%  '

# Row 28
$ws.Range('D28').Value = '0.163'
$ws.Range('E28').Value = '  -1.23%  '

# Row 29
$ws.Range('E29').Value = '  -3.00%  '

# Row 30
$ws.Range('E30').Value = '  +0.08%  '

# Row 31
$ws.Range('E31').Value = '  -2.70%  '

# Row 32
$ws.Range('D32').Value = '0.0₃0859'
$ws.Range('E32').Value = '  +0.16%  '

# Row 33
$ws.Range('E33').Value = '  -1.60%  '

# Row 34
$ws.Range('D34').Value = '5.19'

# Row 35
$ws.Range('D35').Value = '165.30'
$ws.Range('E35').Value = '  -1.91%  '

# Row 36
$ws.Range('D36').Value = '0.414'
$ws.Range('E36').Value = '  +1.21%  '

# Row 37
$ws.Range('D37').Value = '1.00'

# Row 38
$ws.Range('D38').Value = '19.40'
$ws.Range('E38').Value = '  +0.11%  '

# Row 39
$ws.Range('E39').Value = '  -4.09%  '

# Row 40
$ws.Range('E40').Value = '  +0.01%  '

# Row 41
$ws.Range('D41').Value = '165.55'
$ws.Range('E41').Value = '  -0.82%  '

# Row 42
$ws.Range('D42').Value = '39.72'
$ws.Range('E42').Value = '  -1.52%  '

# Row 43
$ws.Range('D43').Value = '3.97'
$ws.Range('E43').Value = '  +2.70%  '

# Row 44
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').Value = '0.0592'
$ws.Range('E44').Value = '  +2.95%  '

# Row 45
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '22.70'
$ws.Range('E45').Value = '  +2.95%  '

# Row 46
$ws.Range('E46').Value = '  -0.53%  '

# Row 47
$ws.Range('E47').Value = '  +2.05%  '

# Row 48
$ws.Range('E48').Value = '  +0.75%  '

# Row 49
$ws.Range('E49').Value = '  -0.66%  '

# Row 50
$ws.Range('E50').Value = '  -0.64%  '

# Row 51
$ws.Range('D51').Value = '0.0₆0233'
$ws.Range('E51').Value = '  +16.05%  '
